$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 829.5454999999999
$ws.Range("J17").Value = 829.5454999999999
$ws.Range("L17").Value = 2488.6365
$ws.Range("N17").Value = -2824.6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10145.956
$ws.Range("J32").Value = 9716.933999999999
$ws.Range("L32").Value = 9716.933999999999
$ws.Range("N32").Value = -10368.934

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 426.6875
$ws.Range("I33").Value = 85.111115
$ws.Range("J33").Value = 865.8570999999999
$ws.Range("K33").Value = 85.111115
$ws.Range("L33").Value = 865.8570999999999
$ws.Range("M33").Value = 143.888885
$ws.Range("N33").Value = -1323.8571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1802.069
$ws.Range("I98").Value = 1652.1428
$ws.Range("J98").Value = 6000
$ws.Range("K98").Value = 1652.1428
$ws.Range("L98").Value = 6000
$ws.Range("M98").Value = -154.1428000000001
$ws.Range("N98").Value = -8996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1364.5714
$ws.Range("I100").Value = 1264.091
$ws.Range("K100").Value = 1264.091
$ws.Range("M100").Value = -723.0909999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 337
$ws.Range("I111").Value = 361.5
$ws.Range("K111").Value = 1084.5
$ws.Range("M111").Value = 1982.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10116.366
$ws.Range("I116").Value = 12961
$ws.Range("K116").Value = 12961
$ws.Range("M116").Value = -9519

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1802.069
$ws.Range("I122").Value = 1652.1428
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 4956.428400000001
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -2506.428400000001
$ws.Range("N122").Value = -22900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 668.0714
$ws.Range("I2").Value = 673.2
$ws.Range("J2").Value = 655.25
$ws.Range("K2").Value = 673.2
$ws.Range("L2").Value = 655.25
$ws.Range("M2").Value = -560.2
$ws.Range("N2").Value = -881.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29460.87
$ws.Range("I32").Value = 30847.584
$ws.Range("J32").Value = 4500
$ws.Range("K32").Value = 30847.584
$ws.Range("L32").Value = 4500
$ws.Range("M32").Value = -30560.584
$ws.Range("N32").Value = -5074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 68251
$ws.Range("J43").Value = 92376.5
$ws.Range("L43").Value = 92376.5
$ws.Range("N43").Value = -93002.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4644.3335
$ws.Range("I45").Value = 3149
$ws.Range("K45").Value = 3149
$ws.Range("M45").Value = -2772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4449.5
$ws.Range("I63").Value = 2447.5
$ws.Range("J63").Value = 4783.1665
$ws.Range("K63").Value = 2447.5
$ws.Range("L63").Value = 4783.1665
$ws.Range("M63").Value = -1761.5
$ws.Range("N63").Value = -6155.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4449.5
$ws.Range("I66").Value = 2447.5
$ws.Range("J66").Value = 4783.1665
$ws.Range("K66").Value = 12237.5
$ws.Range("L66").Value = 23915.8325
$ws.Range("M66").Value = -8805.5
$ws.Range("N66").Value = -30779.8325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2838.125
$ws.Range("I97").Value = 2627.5
$ws.Range("J97").Value = 3048.75
$ws.Range("K97").Value = 2627.5
$ws.Range("L97").Value = 3048.75
$ws.Range("M97").Value = -2131.5
$ws.Range("N97").Value = -4040.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 250000
$ws.Range("J109").Value = 250000
$ws.Range("L109").Value = 250000
$ws.Range("N109").Value = -252774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1099.875
$ws.Range("I110").Value = 1134.3334
$ws.Range("K110").Value = 1134.3334
$ws.Range("M110").Value = 910.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 668.0714
$ws.Range("I116").Value = 673.2
$ws.Range("J116").Value = 655.25
$ws.Range("K116").Value = 673.2
$ws.Range("L116").Value = 655.25
$ws.Range("M116").Value = 1620.8
$ws.Range("N116").Value = -5243.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1562.6586
$ws.Range("I132").Value = 1363.2424
$ws.Range("K132").Value = 4089.7272
$ws.Range("M132").Value = -1559.7272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 668.0714
$ws.Range("I3").Value = 673.2
$ws.Range("J3").Value = 655.25
$ws.Range("K3").Value = 673.2
$ws.Range("L3").Value = 655.25
$ws.Range("M3").Value = -559.2
$ws.Range("N3").Value = -883.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2069.5881
$ws.Range("I99").Value = 2033.75
$ws.Range("J99").Value = 2155.6
$ws.Range("K99").Value = 2033.75
$ws.Range("L99").Value = 2155.6
$ws.Range("M99").Value = -535.75
$ws.Range("N99").Value = -5151.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6074.591
$ws.Range("I107").Value = 6105.3125
$ws.Range("J107").Value = 5992.6665
$ws.Range("K107").Value = 6105.3125
$ws.Range("L107").Value = 5992.6665
$ws.Range("M107").Value = -4185.3125
$ws.Range("N107").Value = -9832.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5884470
$ws.Range("I31").Value = 10000999
$ws.Range("K31").Value = 10000999
$ws.Range("M31").Value = -10000704

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5884470
$ws.Range("I34").Value = 10000999
$ws.Range("K34").Value = 10000999
$ws.Range("M34").Value = -10000797

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13573.63
$ws.Range("J58").Value = 48451.855
$ws.Range("L58").Value = 48451.855
$ws.Range("N58").Value = -48857.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5230.857
$ws.Range("I99").Value = 4171.077
$ws.Range("J99").Value = 6953
$ws.Range("K99").Value = 4171.077
$ws.Range("L99").Value = 6953
$ws.Range("M99").Value = -2673.077
$ws.Range("N99").Value = -9949

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 49244.875
$ws.Range("I105").Value = 55565.57
$ws.Range("K105").Value = 55565.57
$ws.Range("M105").Value = -53818.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1616
$ws.Range("I122").Value = 1600.3684
$ws.Range("J122").Value = 1675.4
$ws.Range("K122").Value = 4801.1052
$ws.Range("L122").Value = 5026.200000000001
$ws.Range("M122").Value = -2351.1052
$ws.Range("N122").Value = -9926.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5230.857
$ws.Range("I126").Value = 4171.077
$ws.Range("J126").Value = 6953
$ws.Range("K126").Value = 12513.231
$ws.Range("L126").Value = 20859
$ws.Range("M126").Value = -10043.231
$ws.Range("N126").Value = -25799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 59888.65
$ws.Range("I132").Value = 63325.188
$ws.Range("K132").Value = 189975.564
$ws.Range("M132").Value = -187445.564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1906.1111
$ws.Range("I134").Value = 1624.5652
$ws.Range("K134").Value = 4873.6956
$ws.Range("M134").Value = -2338.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13573.63
$ws.Range("J136").Value = 48451.855
$ws.Range("L136").Value = 145355.565
$ws.Range("N136").Value = -150455.565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 118989
$ws.Range("J140").Value = 118989
$ws.Range("L140").Value = 118989
$ws.Range("N140").Value = -129349

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4653.0625
$ws.Range("J69").Value = 4696.6
$ws.Range("L69").Value = 14089.8
$ws.Range("N69").Value = -15711.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 4653.0625
$ws.Range("J72").Value = 4696.6
$ws.Range("L72").Value = 42269.4
$ws.Range("N72").Value = -50381.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2910.1667
$ws.Range("I133").Value = 2303.8572
$ws.Range("K133").Value = 6911.571599999999
$ws.Range("M133").Value = -1851.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3017.75
$ws.Range("I140").Value = 3017.75
$ws.Range("K140").Value = 9053.25
$ws.Range("M140").Value = -3873.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -26996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 125000
$ws.Range("N83").Value = -134984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 806.7273
$ws.Range("I97").Value = 831.55554
$ws.Range("K97").Value = 831.55554
$ws.Range("M97").Value = -335.55554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3539.2
$ws.Range("I126").Value = 3056.7144
$ws.Range("J126").Value = 4665
$ws.Range("K126").Value = 9170.143199999999
$ws.Range("L126").Value = 13995
$ws.Range("M126").Value = -6700.143199999999
$ws.Range("N126").Value = -18935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 92174.836
$ws.Range("J133").Value = 93268
$ws.Range("L133").Value = 93268
$ws.Range("N133").Value = -103388

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 106666.336
$ws.Range("J135").Value = 106666.336
$ws.Range("L135").Value = 106666.336
$ws.Range("N135").Value = -116806.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3445.1843
$ws.Range("I136").Value = 3155.8438
$ws.Range("J136").Value = 4988.3335
$ws.Range("K136").Value = 9467.5314
$ws.Range("L136").Value = 14965.0005
$ws.Range("M136").Value = -6917.5314
$ws.Range("N136").Value = -20065.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 993.3333
$ws.Range("I96").Value = 400
$ws.Range("J96").Value = 1290
$ws.Range("K96").Value = 400
$ws.Range("L96").Value = 1290
$ws.Range("M96").Value = 973
$ws.Range("N96").Value = -4036

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1889.2174
$ws.Range("I132").Value = 1505.3572
$ws.Range("K132").Value = 4516.071599999999
$ws.Range("M132").Value = -1986.071599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 23494.68
$ws.Range("I136").Value = 25994
$ws.Range("K136").Value = 77982
$ws.Range("M136").Value = -75432
